# Adds a new "2022-Q3" quarterly sheet (with its fund-holdings detail) and
# records it as the newest row in the "总计" (totals) summary sheet, pushing
# every existing quarter down by one row / one tab position.

function Set-TextCell($ws, $addr, $val) {
    # Forces a numeric-looking string (e.g. "44.76") to be stored as TEXT
    # instead of being auto-coerced into a number by Excel's input parser.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new row 2 for 2022-Q3, shifting the
#    existing 2022-Q2 .. 2020-Q4 rows down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 1.84

# Give the new A2 the same centred/bordered look as the rest of column A.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Column A is a plain 0-based row counter; Insert() left the old values in
# place (0,1,2,...) on rows 3-9 instead of bumping them, so renumber them.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("A8").Value = 6
$totalSheet.Range("A9").Value = 7

# ---------------------------------------------------------------------
# 2) New "2022-Q3" fund-holdings sheet, inserted right after "总计" (i.e.
#    right before the existing "2022-Q2" tab). Duplicate "2022-Q2" so the
#    new sheet starts with identical layout/styling, then overwrite the
#    data with the 2022-Q3 numbers.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The template only has 5 data rows (rows 2-6); 2022-Q3 needs 9 (rows 2-10).
# Stretch it by copying row 6's formatting down to the new rows 7-10.
$newSheet.Range("A6:H6").Copy()
$newSheet.Range("A7:H10").PasteSpecial(-4122)

$rows = @(
  @(0, "512980", "广发中证传媒ETF", "44.76", "99.29", "3.34", "1.4950", 5),
  @(1, "160629", "鹏华中证传媒指数（LOF）A", "6.41", "94.58", "3.15", "0.2019", 5),
  @(2, "159805", "鹏华中证传媒ETF", "1.71", "98.37", "3.30", "0.0564", 5),
  @(3, "164818", "工银瑞信中证传媒指数（LOF）A", "1.65", "93.46", "3.12", "0.0515", 5),
  @(4, "159725", "工银瑞信中证线上消费主题ETF", "0.57", "98.42", "2.44", "0.0139", 10),
  @(5, "517770", "浦银安盛中证沪港深游戏及文化传媒ETF", "0.32", "93.49", "2.21", "0.0071", 10),
  @(6, "010677", "工银瑞信中证传媒指数（LOF）C", "0.21", "93.46", "3.12", "0.0066", 5),
  @(7, "015675", "鹏华中证传媒指数（LOF）C", "0.17", "94.58", "3.15", "0.0054", 5),
  @(8, "516190", "华夏中证文娱传媒ETF", "0.13", "96.01", "2.72", "0.0035", 8)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row[0]
    Set-TextCell $newSheet "B$r" $row[1]
    Set-TextCell $newSheet "C$r" $row[2]
    Set-TextCell $newSheet "D$r" $row[3]
    Set-TextCell $newSheet "E$r" $row[4]
    Set-TextCell $newSheet "F$r" $row[5]
    Set-TextCell $newSheet "G$r" $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r++
}

# Keep "总计" as the active/selected sheet, same as before the edit.
$totalSheet.Activate()
$totalSheet.Range("A1").Select()
